$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.147.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.909.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3890"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07859"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9918"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.049"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.755"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009940"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.154.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.327"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.142.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.085"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.933"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.877"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09323"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8964"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.226"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.44%  "

$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.146"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05787"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.168"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02089"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5711"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.672"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1811"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.733"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002865"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +73.57%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.98%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5353"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.170"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.18%  "

$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.842"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.551"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
